$d = $word.ActiveDocument

# ------------------------------------------------------------------
# Locate the two anchor paragraphs dynamically (robust to any minor
# index drift): the "Group ... Presentation slides ... Done" list item
# (numId=2) under "Later work:" and the "Final work:" paragraph that
# currently follows the four placeholder "???" list items.
# ------------------------------------------------------------------
$groupParaIndex = -1
$finalWorkParaIndex = -1
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $p = $d.Paragraphs.Item($i)
    $t = $p.Range.Text
    if ($groupParaIndex -eq -1 -and $t -like "Group*Presentation slides*Done*") {
        $groupParaIndex = $i
    }
    if ($finalWorkParaIndex -eq -1 -and $groupParaIndex -ne -1 -and $t -like "Final work:*") {
        $finalWorkParaIndex = $i
    }
}

if ($groupParaIndex -eq -1 -or $finalWorkParaIndex -eq -1) {
    throw "Could not locate anchor paragraphs (group=$groupParaIndex, finalWork=$finalWorkParaIndex)"
}

# ------------------------------------------------------------------
# Delete the four placeholder list paragraphs ("??? ... Hand outs",
# "??? ... video", "??? ... whiteboard content", "??? ... ???") that
# sit between the "Group" paragraph and the "Final work:" paragraph.
# ------------------------------------------------------------------
$firstDeleteIndex = $groupParaIndex + 1
$lastDeleteIndex = $finalWorkParaIndex - 1

if ($lastDeleteIndex -ge $firstDeleteIndex) {
    $startRange = $d.Paragraphs.Item($firstDeleteIndex).Range.Start
    $endRange = $d.Paragraphs.Item($lastDeleteIndex).Range.End
    $deleteRange = $d.Range($startRange, $endRange)
    $deleteRange.Delete()
}

# ------------------------------------------------------------------
# Move the (hidden) "_GoBack" bookmark from right after "Group" to the
# very start of the "Final work:" paragraph. Re-adding a bookmark with
# the same name relocates it (a document can only have one bookmark
# per name), which also removes it from its old position. Re-locate
# the "Final work:" paragraph by searching again, since paragraph /
# range handles are index-based and do not track live edits made by
# the deletion above.
# ------------------------------------------------------------------
$finalWorkParaIndex2 = -1
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $t = $d.Paragraphs.Item($i).Range.Text
    if ($t -like "Final work:*") {
        $finalWorkParaIndex2 = $i
        break
    }
}
if ($finalWorkParaIndex2 -eq -1) {
    throw "Could not relocate the 'Final work:' paragraph after deletion"
}

$finalWorkPara = $d.Paragraphs.Item($finalWorkParaIndex2)
$bookmarkRange = $d.Range($finalWorkPara.Range.Start, $finalWorkPara.Range.Start)
$d.Bookmarks.Add("_GoBack", $bookmarkRange) | Out-Null
